$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New audio entry row (row 46): Menu UI click sound
$ws.Range("A46").Value = "448081__breviceps__tic-toc-click.wav"
$ws.Range("B46").Value = "Menu UI Click"
$ws.Range("C46").Value = "https://freesound.org/"
$ws.Range("D46").Value = "https://freesound.org/people/Breviceps/sounds/448081/"
$ws.Range("E46").Value = "https://creativecommons.org/publicdomain/zero/1.0/"

# Match the styling already used by the other rows in the table
$ws.Range("A46:B46").Style = $ws.Range("A45:B45").Style
$ws.Range("C46:E46").Style = $ws.Range("C45:E45").Style

# Hyperlinks for the website / direct link / license cells
$ws.Hyperlinks.Add($ws.Range("C46"), "https://freesound.org/")
$ws.Hyperlinks.Add($ws.Range("E46"), "https://creativecommons.org/publicdomain/zero/1.0/")
$ws.Hyperlinks.Add($ws.Range("D46"), "https://freesound.org/people/Breviceps/sounds/448081/")

# Restore table cell styling that the hyperlink insertion may have touched
$ws.Range("C46:E46").Style = $ws.Range("C45:E45").Style

# Scroll/selection state recorded for the sheet view
$ws.Application.ActiveWindow.ScrollRow = 35
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D47").Select()
